$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster table lists one player per row (row 15 = "John Butler (TW)",
# row 16 = "Kevin Knox"). The edit re-sorts those two players, so every
# attribute in columns B:K (No., Player, Pos, Ht, Wt, Birth Date, the
# unnamed country column, Exp, College, bbref url) needs to trade places
# between row 15 and row 16; column A (the row's positional index) stays
# put.
#
# Using Range.Value to move the data would make Excel re-infer each cell's
# type, turning text that merely looks numeric (e.g. the Exp column's "4")
# into a real number. Range.Copy moves the literal cell contents (value +
# type) without that re-inference, so it is used here together with an
# off-table scratch row to perform a true swap.

$row15 = $ws.Range("B15:K15")
$row16 = $ws.Range("B16:K16")
$scratch = $ws.Range("B100:K100")

$row15.Copy($scratch)
$row16.Copy($row15)
$scratch.Copy($row16)
$scratch.Clear()
